$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34, pushing the existing rows 34-39 down to 35-40.
$ws.Rows.Item(34).Insert()

# Populate the new row 34 with the weekly record (same fixed attributes as the
# surrounding Chirimoya / Terminal Hortofrutícola Agro Chillán rows).
$ws.Range("A34").Value = 7
$ws.Range("B34").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C34").Value = "Ñuble"
$ws.Range("D34").Value = 45244
$ws.Range("E34").Value = 16
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100107
$ws.Range("H34").Value = "Otros"
$ws.Range("I34").Value = 100107002
$ws.Range("J34").Value = "Chirimoya"
$ws.Range("K34").Value = "Cultivar IV Región"
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 90
$ws.Range("N34").Value = 20000
$ws.Range("O34").Value = 21000
$ws.Range("P34").Value = 20444
$ws.Range("Q34").Value = "$/bandeja 10 kilos"
$ws.Range("R34").Value = "Provincia de Limarí"
$ws.Range("S34").Value = 2044
$ws.Range("T34").Value = 10

# Match the date-style formatting used by the rest of column D.
$ws.Range("D34").NumberFormat = $ws.Range("D35").NumberFormat
